# Update Betfair Back/Lay odds for the games of 2025-12-17
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Young Boys vs Grasshoppers Zurich
$ws.Range("F2").Value = 1.61
$ws.Range("G2").Value = 1.63
$ws.Range("H2").Value = 5.1
$ws.Range("J2").Value = 4.6
$ws.Range("Q2").Value = 1.5

# Row 3 - FC Zurich vs Lugano
$ws.Range("G3").Value = 3.3
$ws.Range("I3").Value = 2.86
$ws.Range("J3").Value = 3.55
$ws.Range("P3").Value = 2.22
$ws.Range("Q3").Value = 1.65

# Row 4 - Luzern vs FC Basel
$ws.Range("I4").Value = 2.28
$ws.Range("P4").Value = 2.66

# Row 5 - Dundee Utd vs Celtic
$ws.Range("F5").Value = 6.8
$ws.Range("G5").Value = 8
$ws.Range("H5").Value = 1.49
$ws.Range("I5").Value = 1.52
$ws.Range("J5").Value = 5
$ws.Range("K5").Value = 5.4
$ws.Range("O5").Value = 1.2
$ws.Range("T5").Value = 1.83
$ws.Range("U5").Value = 2.04
$ws.Range("X5").Value = 28
$ws.Range("Y5").Value = 12.5
$ws.Range("AB5").Value = 32
$ws.Range("AC5").Value = 14
$ws.Range("AD5").Value = 12.5
$ws.Range("AE5").Value = 17.5
$ws.Range("AG5").Value = 28
$ws.Range("AI5").Value = 36
$ws.Range("AL5").Value = 90
$ws.Range("AM5").Value = 130
